$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.288150666666667
$ws.Range("H2").Value = 6.864452
$ws.Range("I2").Value = 0.3964219041944151
$ws.Range("J2").Value = 0.3964219041944151
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 150.1098686666667
$ws.Range("N2").Value = 450.329606
$ws.Range("O2").Value = 0.7276622610660995
$ws.Range("P2").Value = 0.7276622610660997
$ws.Range("Q2").Value = 343.4739960628791
$ws.Range("R2").Value = 3091.265964565912
$ws.Range("S2").Value = 0.2884612591422367
$ws.Range("T2").Value = 0.2884612591422368

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.288150666666667
$ws.Range("H3").Value = 6.864452
$ws.Range("I3").Value = 0.3964219041944151
$ws.Range("J3").Value = 0.3964219041944151
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 42.32476666666667
$ws.Range("N3").Value = 126.9743
$ws.Range("O3").Value = 0.2051706239258123
$ws.Range("P3").Value = 0.2051706239258124
$ws.Range("Q3").Value = 96.84544306484445
$ws.Range("R3").Value = 871.6089875836
$ws.Range("S3").Value = 0.08133412942142675
$ws.Range("T3").Value = 0.08133412942142676

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.288150666666667
$ws.Range("H4").Value = 6.864452
$ws.Range("I4").Value = 0.3964219041944151
$ws.Range("J4").Value = 0.3964219041944151
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.940565666666667
$ws.Range("N4").Value = 14.821697
$ws.Range("O4").Value = 0.02394954586187395
$ws.Range("P4").Value = 0.02394954586187395
$ws.Range("Q4").Value = 11.30475862389378
$ws.Range("R4").Value = 101.742827615044
$ws.Range("S4").Value = 0.009494124575155544
$ws.Range("T4").Value = 0.009494124575155546

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.288150666666667
$ws.Range("H5").Value = 6.864452
$ws.Range("I5").Value = 0.3964219041944151
$ws.Range("J5").Value = 0.3964219041944151
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.915377333333334
$ws.Range("N5").Value = 26.746132
$ws.Range("O5").Value = 0.04321756914621411
$ws.Range("P5").Value = 0.04321756914621412
$ws.Range("Q5").Value = 20.39972658885156
$ws.Range("R5").Value = 183.597539299664
$ws.Range("S5").Value = 0.017132391055596
$ws.Range("T5").Value = 0.017132391055596

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.588894
$ws.Range("H6").Value = 7.766681999999999
$ws.Range("I6").Value = 0.4485256605643812
$ws.Range("J6").Value = 0.4485256605643813
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 150.1098686666667
$ws.Range("N6").Value = 450.329606
$ws.Range("O6").Value = 0.7276622610660995
$ws.Range("P6").Value = 0.7276622610660997
$ws.Range("Q6").Value = 388.6185383319213
$ws.Range("R6").Value = 3497.566844987292
$ws.Range("S6").Value = 0.3263751963124435
$ws.Range("T6").Value = 0.3263751963124436

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.588894
$ws.Range("H7").Value = 7.766681999999999
$ws.Range("I7").Value = 0.4485256605643812
$ws.Range("J7").Value = 0.4485256605643813
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 42.32476666666667
$ws.Range("N7").Value = 126.9743
$ws.Range("O7").Value = 0.2051706239258123
$ws.Range("P7").Value = 0.2051706239258124
$ws.Range("Q7").Value = 109.5743344747333
$ws.Range("R7").Value = 986.1690102726
$ws.Range("S7").Value = 0.09202428962473122
$ws.Range("T7").Value = 0.09202428962473123

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.588894
$ws.Range("H8").Value = 7.766681999999999
$ws.Range("I8").Value = 0.4485256605643812
$ws.Range("J8").Value = 0.4485256605643813
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.940565666666667
$ws.Range("N8").Value = 14.821697
$ws.Range("O8").Value = 0.02394954586187395
$ws.Range("P8").Value = 0.02394954586187395
$ws.Range("Q8").Value = 12.79060081103933
$ws.Range("R8").Value = 115.115407299354
$ws.Range("S8").Value = 0.01074198587791395
$ws.Range("T8").Value = 0.01074198587791396

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.588894
$ws.Range("H9").Value = 7.766681999999999
$ws.Range("I9").Value = 0.4485256605643812
$ws.Range("J9").Value = 0.4485256605643813
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.915377333333334
$ws.Range("N9").Value = 26.746132
$ws.Range("O9").Value = 0.04321756914621411
$ws.Range("P9").Value = 0.04321756914621412
$ws.Range("Q9").Value = 23.08096688600267
$ws.Range("R9").Value = 207.728701974024
$ws.Range("S9").Value = 0.01938418874929251
$ws.Range("T9").Value = 0.01938418874929251

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3123523333333333
$ws.Range("H10").Value = 0.9370569999999999
$ws.Range("I10").Value = 0.05411501461132016
$ws.Range("J10").Value = 0.05411501461132018
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 150.1098686666667
$ws.Range("N10").Value = 450.329606
$ws.Range("O10").Value = 0.7276622610660995
$ws.Range("P10").Value = 0.7276622610660997
$ws.Range("Q10").Value = 46.88716773439355
$ws.Range("R10").Value = 421.984509609542
$ws.Range("S10").Value = 0.03937745388969825
$ws.Range("T10").Value = 0.03937745388969826

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3123523333333333
$ws.Range("H11").Value = 0.9370569999999999
$ws.Range("I11").Value = 0.05411501461132016
$ws.Range("J11").Value = 0.05411501461132018
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 42.32476666666667
$ws.Range("N11").Value = 126.9743
$ws.Range("O11").Value = 0.2051706239258123
$ws.Range("P11").Value = 0.2051706239258124
$ws.Range("Q11").Value = 13.22023962612222
$ws.Range("R11").Value = 118.9821566351
$ws.Range("S11").Value = 0.01110281131155901
$ws.Range("T11").Value = 0.01110281131155901

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.3123523333333333
$ws.Range("H12").Value = 0.9370569999999999
$ws.Range("I12").Value = 0.05411501461132016
$ws.Range("J12").Value = 0.05411501461132018
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.940565666666667
$ws.Range("N12").Value = 14.821697
$ws.Range("O12").Value = 0.02394954586187395
$ws.Range("P12").Value = 0.02394954586187395
$ws.Range("Q12").Value = 1.543197213969889
$ws.Range("R12").Value = 13.888774925729
$ws.Range("S12").Value = 0.001296030024249791
$ws.Range("T12").Value = 0.001296030024249792

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.3123523333333333
$ws.Range("H13").Value = 0.9370569999999999
$ws.Range("I13").Value = 0.05411501461132016
$ws.Range("J13").Value = 0.05411501461132018
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.915377333333334
$ws.Range("N13").Value = 26.746132
$ws.Range("O13").Value = 0.04321756914621411
$ws.Range("P13").Value = 0.04321756914621412
$ws.Range("Q13").Value = 2.784738912613777
$ws.Range("R13").Value = 25.062650213524
$ws.Range("S13").Value = 0.002338719385813116
$ws.Range("T13").Value = 0.002338719385813117

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5826116666666666
$ws.Range("H14").Value = 1.747835
$ws.Range("I14").Value = 0.1009374206298835
$ws.Range("J14").Value = 0.1009374206298836
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 150.1098686666667
$ws.Range("N14").Value = 450.329606
$ws.Range("O14").Value = 0.7276622610660995
$ws.Range("P14").Value = 0.7276622610660997
$ws.Range("Q14").Value = 87.45576076700111
$ws.Range("R14").Value = 787.1018469030099
$ws.Range("S14").Value = 0.07344835172172103
$ws.Range("T14").Value = 0.07344835172172104

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5826116666666666
$ws.Range("H15").Value = 1.747835
$ws.Range("I15").Value = 0.1009374206298835
$ws.Range("J15").Value = 0.1009374206298836
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 42.32476666666667
$ws.Range("N15").Value = 126.9743
$ws.Range("O15").Value = 0.2051706239258123
$ws.Range("P15").Value = 0.2051706239258124
$ws.Range("Q15").Value = 24.65890284894444
$ws.Range("R15").Value = 221.9301256405
$ws.Range("S15").Value = 0.02070939356809537
$ws.Range("T15").Value = 0.02070939356809538

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5826116666666666
$ws.Range("H16").Value = 1.747835
$ws.Range("I16").Value = 0.1009374206298835
$ws.Range("J16").Value = 0.1009374206298836
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 4.940565666666667
$ws.Range("N16").Value = 14.821697
$ws.Range("O16").Value = 0.02394954586187395
$ws.Range("P16").Value = 0.02394954586187395
$ws.Range("Q16").Value = 2.878431197332778
$ws.Range("R16").Value = 25.905880775995
$ws.Range("S16").Value = 0.002417405384554657
$ws.Range("T16").Value = 0.002417405384554658

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5826116666666666
$ws.Range("H17").Value = 1.747835
$ws.Range("I17").Value = 0.1009374206298835
$ws.Range("J17").Value = 0.1009374206298836
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.915377333333334
$ws.Range("N17").Value = 26.746132
$ws.Range("O17").Value = 0.04321756914621411
$ws.Range("P17").Value = 0.04321756914621412
$ws.Range("Q17").Value = 5.194202847135555
$ws.Range("R17").Value = 46.74782562421999
$ws.Range("S17").Value = 0.004362269955512491
$ws.Range("T17").Value = 0.004362269955512492
